# Checkliste Architektur: add "Status" column D with check-off values for
# the first architecture review (with Jochen, 10.06.2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain "ja" confirmations for most checklist rows.
$okRows = 4,7,8,9,11,12,13,14,15,16,17,18
foreach ($r in $okRows) {
    $ws.Cells.Item($r, 4).Value = "ja"
}

# Row 5: accepted with a caveat.
$ws.Range("D5").Value = "ja, angepasst"

# Rows 6 and 10 are section headers (bold question) - mark with bold "ja".
$ws.Range("D6").Value = "ja"
$ws.Range("D6").Font.Bold = $true
$ws.Range("D10").Value = "ja"
$ws.Range("D10").Font.Bold = $true

# Row 19: this checklist item will be dropped going forward.
$ws.Range("D19").Value = "wird gestrichen"

# Row 21: audit note for this pass.
$ws.Range("D21").Value = "gecheckt am 10.06.2016"

# Size column D to fit its new contents, matching cols B and C (best-fit
# width for the longest entry, "gecheckt am 10.06.2016").
$ws.Columns("D").ColumnWidth = 21

# Match the print setup recorded for this check.
$ws.PageSetup.PaperSize = 77
$ws.PageSetup.Orientation = 1

# Leave the cursor where the reviewer left off.
[void]$ws.Range("D18").Select()
